$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.9916784696904799
$ws.Range("C1").Value = 0.92306296516500552
$ws.Range("AO1").Value = 0.67472899813166598
$ws.Range("BH1").Value = 0.94982907705134934
$ws.Range("BO1").Value = 0.82741154120043992
$ws.Range("C2").Value = 0.73715749211083914
$ws.Range("BP2").Value = 0.9064194316448565
$ws.Range("D3").Value = 0.89386018788839172
$ws.Range("E4").Value = 0.92396526484034647
$ws.Range("C5").Value = 0.74566251532275674
$ws.Range("G5").Value = 0.85986825392762101
$ws.Range("G6").Value = 0.99173540566027363
$ws.Range("F8").Value = 0.88372305311146437
$ws.Range("I8").Value = 0.90864773835326962
$ws.Range("D9").Value = 0.98314690641517077
$ws.Range("G9").Value = 0.96875718355781038
$ws.Range("J9").Value = 0.63286319095431431
$ws.Range("H10").Value = 0.83391131009560004
$ws.Range("L10").Value = 0.97039405114125077
$ws.Range("AE10").Value = 0.87347630490221606
$ws.Range("AY10").Value = 0.82496678228753217
$ws.Range("M11").Value = 0.96134378009796562
$ws.Range("K12").Value = 0.80637262706220136
$ws.Range("N13").Value = 0.66425282491361526
$ws.Range("O13").Value = 0.94847289500427356
$ws.Range("L14").Value = 0.93007647818245198
$ws.Range("BL14").Value = 0.72854609891106548
$ws.Range("N15").Value = 0.76856648195424992
$ws.Range("Q15").Value = 0.71299581805219614
$ws.Range("Q16").Value = 0.85625429859124136
$ws.Range("R17").Value = 0.69530407806459216
$ws.Range("BE17").Value = 0.65860976047042641
$ws.Range("E18").Value = 0.9225154940173208
$ws.Range("P18").Value = 0.69567415309116531
$ws.Range("R19").Value = 0.52721828090056611
$ws.Range("R20").Value = 0.95146438834532676
$ws.Range("S20").Value = 0.98738430215560302
$ws.Range("BL20").Value = 0.90005453483007547
$ws.Range("S21").Value = 0.67913692116533153
$ws.Range("V21").Value = 0.99345350380850928
$ws.Range("AH21").Value = 0.99510311332451273
$ws.Range("BA21").Value = 0.5532619005267102
$ws.Range("T22").Value = 0.96300176969609352
$ws.Range("W22").Value = 0.9241250270135607
$ws.Range("X22").Value = 0.99605717129692839
$ws.Range("I23").Value = 0.98219462911211342
$ws.Range("X23").Value = 0.79007330338830173
$ws.Range("Y23").Value = 0.98365159083100906
$ws.Range("N24").Value = 0.67809584829664316
$ws.Range("AZ24").Value = 0.84217205296733399
$ws.Range("P25").Value = 0.79953458721073378
$ws.Range("Y27").Value = 0.96654930014283513
$ws.Range("AB27").Value = 0.84648657672692829
$ws.Range("L28").Value = 0.5789422295148039
$ws.Range("Z28").Value = 0.69554350734822745
$ws.Range("AC28").Value = 0.84923962880780213
$ws.Range("AD28").Value = 0.79020684499877047
$ws.Range("AA29").Value = 0.71171658073727806
$ws.Range("AC30").Value = 0.76918720779400429
$ws.Range("AE30").Value = 0.82496178407437948
$ws.Range("AC31").Value = 0.81642555926790772
$ws.Range("AF31").Value = 0.99580694900952327
$ws.Range("AV31").Value = 0.72156588540349864
$ws.Range("AD32").Value = 0.72436620930696449
$ws.Range("AG32").Value = 0.86259234723589506
$ws.Range("X33").Value = 0.62026949387535057
$ws.Range("AG34").Value = 0.8775222231494948
$ws.Range("AI34").Value = 0.94749729672077354
$ws.Range("AJ34").Value = 0.89372721199239324
$ws.Range("W35").Value = 0.89883967437136358
$ws.Range("AG35").Value = 0.8342300012263999
$ws.Range("AK35").Value = 0.64918703389416055
$ws.Range("AL36").Value = 0.75846109828724184
$ws.Range("AJ37").Value = 0.59405063933787794
$ws.Range("AF38").Value = 0.94084234510367226
$ws.Range("AK38").Value = 0.70265063861364063
$ws.Range("AM38").Value = 0.98107934479330894
$ws.Range("T39").Value = 0.66836203448002984
$ws.Range("AK39").Value = 0.83171588058708179
$ws.Range("AN39").Value = 0.94130345941364957
$ws.Range("AL40").Value = 0.9792768085937178
$ws.Range("AO40").Value = 0.58212332035422221
$ws.Range("AP40").Value = 0.81485365217605854
$ws.Range("Q41").Value = 0.97642044345745993
$ws.Range("AR42").Value = 0.94074518181423383
$ws.Range("AP43").Value = 0.95128622980250666
$ws.Range("AT44").Value = 0.82931488738059622
$ws.Range("AQ45").Value = 0.69491771597030549
$ws.Range("AR45").Value = 0.95765726483364233
$ws.Range("G46").Value = 0.97305718608398628
$ws.Range("AQ46").Value = 0.72954340665626938
$ws.Range("AS46").Value = 0.92638926502693786
$ws.Range("AU46").Value = 0.99387179906031542
$ws.Range("AV46").Value = 0.83115888291468942
$ws.Range("U47").Value = 0.84528745635314517
$ws.Range("AW47").Value = 0.72871478239266196
$ws.Range("B49").Value = 0.79022893534344418
$ws.Range("AZ49").Value = 0.64612327988376228
$ws.Range("F50").Value = 0.75707904034075679
$ws.Range("AV50").Value = 0.89318875385006447
$ws.Range("AZ50").Value = 0.89532702296934907
$ws.Range("Z51").Value = 0.87166663330298177
$ws.Range("AW51").Value = 0.94348176883402357
$ws.Range("BA51").Value = 0.96786704243900368
$ws.Range("BD51").Value = 0.88396053654512508
$ws.Range("AS52").Value = 0.83794909622267144
$ws.Range("AY52").Value = 0.97295097157323773
$ws.Range("BK52").Value = 0.82452037199215367
$ws.Range("AA53").Value = 0.81290112971571238
$ws.Range("AZ54").Value = 0.96548549015776508
$ws.Range("BC54").Value = 0.95589754817616424
$ws.Range("F55").Value = 0.95403272437392106
$ws.Range("BB56").Value = 0.97260764242412223
$ws.Range("BC56").Value = 0.79638427471371043
$ws.Range("G57").Value = 0.97405721015702973
$ws.Range("O57").Value = 0.86484089954069265
$ws.Range("BD57").Value = 0.91128931394088641
$ws.Range("BG58").Value = 0.93058161936071082
$ws.Range("BE59").Value = 0.94305913055886159
$ws.Range("BH59").Value = 0.99761305866952799
$ws.Range("BI59").Value = 0.91396606077854714
$ws.Range("BF60").Value = 0.97212380543598564
$ws.Range("BJ60").Value = 0.98717148665759447
$ws.Range("BJ61").Value = 0.99750305195943434
$ws.Range("AW63").Value = 0.68988467649630647
$ws.Range("BI63").Value = 0.93586624769931159
$ws.Range("BJ63").Value = 0.7056536115491242
$ws.Range("BB64").Value = 0.96650871678717909
$ws.Range("BJ64").Value = 0.74464259294560808
$ws.Range("BN64").Value = 0.8738720464716152
$ws.Range("D65").Value = 0.75315864371914332
$ws.Range("BD67").Value = 0.95049120192316772
$ws.Range("BM67").Value = 0.93176254784514323
$ws.Range("BN67").Value = 0.86360052104405893
$ws.Range("AO68").Value = 0.65170380713830522
$ws.Range("BN68").Value = 0.93460929354209932
